$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season record columns (AC, AD, AE)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header formatting (bold, border, centered) from an existing
# header cell (AB1) onto the new header cells
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$lastRow = 44

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 103   # AC
    $ws.Cells.Item($r, 30).Value = 59    # AD
    $ws.Cells.Item($r, 31).Value = 0     # AE
}
